# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (column G) values for the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 9
